$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Debug")

# Replace the generic A/B/C/D multiple-choice placeholders with real answer text.
# Row 2 (Who?)   -> Yes / No / Maybe / So
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Range("D2").Value = "Maybe"
$ws.Range("E2").Value = "So"

# Row 3 (What?)  -> So / Yes / No / Maybe
$ws.Range("B3").Value = "So"
$ws.Range("C3").Value = "Yes"
$ws.Range("D3").Value = "No"
$ws.Range("E3").Value = "Maybe"

# Row 4 (When?)  -> Maybe / So / Yes / No
$ws.Range("B4").Value = "Maybe"
$ws.Range("C4").Value = "So"
$ws.Range("D4").Value = "Yes"
$ws.Range("E4").Value = "No"

# Row 5 (Where?) -> No / Maybe / So / Yes
$ws.Range("B5").Value = "No"
$ws.Range("C5").Value = "Maybe"
$ws.Range("D5").Value = "So"
$ws.Range("E5").Value = "Yes"

# Row 6 (Why?)   -> Someone / Else / Is / Ugly
$ws.Range("B6").Value = "Someone"
$ws.Range("C6").Value = "Else"
$ws.Range("D6").Value = "Is"
$ws.Range("E6").Value = "Ugly"

# Row 7 (How?)   -> Ugly / Someone / Else / Is
$ws.Range("B7").Value = "Ugly"
$ws.Range("C7").Value = "Someone"
$ws.Range("D7").Value = "Else"
$ws.Range("E7").Value = "Is"

# Leave the cursor where the author last edited.
$ws.Activate() | Out-Null
$ws.Range("E7").Select() | Out-Null
